$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "AutoOppProbBlank"
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = "New Business`t"

# F4 needs a brand-new style: numFmtId 0, a new font (Aptos Narrow, black, size 11)
$ws.Range("F4").Value = "Call customer"
$ws.Range("F4").Font.Color = 0

$ws.Range("G4").Value = "Prospecting"
$ws.Range("I4").Value = "Regina"

# ---------------------------------------------------------------------------
# Row 5  (row height 29.25)
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 6000
$ws.Range("C5").Value = "Existing"
$ws.Range("F5").Value = "Schedule meeting"
$ws.Range("G5").Value = "qualification"

# I5 needs a brand-new style: numFmtId 49 (text) + wrapText
$ws.Range("I5").Value = "LeadMandatory"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").WrapText = $true

$ws.Range("J5").Value = "LeadMandatory"
$ws.Rows(5).RowHeight = 29.25

# ---------------------------------------------------------------------------
# Row 6  (row height 29.25)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "AutoOppmnt"
$ws.Range("C6").Value = "Existing"
$ws.Range("D6").Value = "Joe Biden"

$ws.Range("E2").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 45936

$ws.Range("F6").Value = "Schedule meeting"
$ws.Range("G6").Value = "qualification"

# H6 holds "30" as TEXT (column H default style is already numFmtId 49 => text)
$ws.Range("H6").Value = "30"

# I6 reuses the new wrap/text style created at I5
$ws.Range("I5").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Value = "LeadMandatory"

# J6 reuses the new font style created at F4
$ws.Range("J6").Value = "Renewal for key account"
$ws.Range("F4").Copy()
$ws.Range("J6").PasteSpecial(-4122)

$ws.Rows(6).RowHeight = 29.25

# ---------------------------------------------------------------------------
# Row 7  (row height 29.25)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "AutoOpptemp"
$ws.Range("F4").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("B7").Value = 10000
$ws.Range("D7").Value = "Mahesh Patel"

$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 45879

$ws.Range("F7").Value = "Schedule meeting"
$ws.Range("G7").Value = "qualification"
$ws.Range("H7").Value = "31"

$ws.Range("I5").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = "LeadMandatory"

$ws.Range("J7").Value = "Renewal for key account"
$ws.Range("F4").Copy()
$ws.Range("J7").PasteSpecial(-4122)

$ws.Rows(7).RowHeight = 29.25

# ---------------------------------------------------------------------------
# Row 8
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "AutoOpptemper"
$ws.Range("B8").Value = 2000
$ws.Range("C8").Value = "Existing"
$ws.Range("D8").Value = "Lily"

$ws.Range("E2").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = 46013

$ws.Range("G8").Value = "qualification"

# H8 must stay numeric (style text-format column, but value is a real number)
$ws.Range("H8").Style = "Normal"
$ws.Range("H8").Value = 30
$ws.Range("H2").Copy()
$ws.Range("H8").PasteSpecial(-4122)

# I8 reuses the existing Arial style (same as I2)
$ws.Range("I2").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = "Regina"

$ws.Range("J8").Value = "New opportunity from campaign"

# ---------------------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "AutoOppfifa"
$ws.Range("B9").Value = 3000
$ws.Range("C9").Value = "New Business"
$ws.Range("D9").Value = "sumit"

$ws.Range("F9").Value = "Call customer"
$ws.Range("F4").Copy()
$ws.Range("F9").PasteSpecial(-4122)

$ws.Range("H9").Value = "45"

$ws.Range("I2").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = "Regina"

$ws.Range("J9").Value = "New opportunity from campaign"

# ---------------------------------------------------------------------------
# Row 10
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "AutoOppTst"
$ws.Range("B10").Value = 89990
$ws.Range("C10").Value = "New Business"
$ws.Range("D10").Value = "sumit"

$ws.Range("E2").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = 45886

$ws.Range("F10").Value = "Schedule meeting"
$ws.Range("G10").Value = "Prospecting"
$ws.Range("H10").Value = "10"

$ws.Range("J10").Value = "Renewal for key account"
$ws.Range("F4").Copy()
$ws.Range("J10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Column widths (closest reachable values given the engine's 1/6 quantization)
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 16.5
$ws.Columns("C").ColumnWidth = 17.166666666666668
$ws.Columns("D").ColumnWidth = 12.0

# ---------------------------------------------------------------------------
# Final selection (drives dimension + selection in the sheet view)
# ---------------------------------------------------------------------------
$ws.Range("J10").Select()
